# Update Facebook and Twitter pivot-table data per commit "updated Facebook and Twitter data"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("C4").Value = 39647
$ws.Range("D4").Value = 3604.3
$ws.Range("E4").Value = 3706.5
$ws.Range("G4").Value = 638.5
$ws.Range("H4").Value = 3054
$ws.Range("I4").Value = 4806
$ws.Range("J4").Value = 12666
$ws.Range("L4").Value = 4955.9
$ws.Range("P4").Value = 70240
$ws.Range("Q4").Value = 6385.5
$ws.Range("R4").Value = 6737.7
$ws.Range("T4").Value = 531.5
$ws.Range("U4").Value = 4132
$ws.Range("V4").Value = 10066.5
$ws.Range("W4").Value = 20340
$ws.Range("Y4").Value = 7804.4
$ws.Range("AC4").Value = 13374
$ws.Range("AD4").Value = 1215.8
$ws.Range("AE4").Value = 1092.8
$ws.Range("AG4").Value = 196
$ws.Range("AH4").Value = 1201
$ws.Range("AI4").Value = 1873
$ws.Range("AJ4").Value = 3465
$ws.Range("AL4").Value = 1671.8
$ws.Range("AO4").Value = 0.2

# Row 5
$ws.Range("C5").Value = 718118
$ws.Range("D5").Value = 1841.3
$ws.Range("E5").Value = 2721.3
$ws.Range("G5").Value = 268.2
$ws.Range("H5").Value = 1010
$ws.Range("I5").Value = 2350.8
$ws.Range("J5").Value = 23112
$ws.Range("L5").Value = 2258.2
$ws.Range("M5").Value = 318
$ws.Range("N5").Value = 81.5
$ws.Range("P5").Value = 271508
$ws.Range("Q5").Value = 696.2
$ws.Range("R5").Value = 1473.6
$ws.Range("T5").Value = 50.5
$ws.Range("U5").Value = 238.5
$ws.Range("V5").Value = 692.8
$ws.Range("W5").Value = 13652
$ws.Range("Y5").Value = 838
$ws.Range("Z5").Value = 324
$ws.Range("AA5").Value = 83.09999999999999
$ws.Range("AC5").Value = 314407
$ws.Range("AD5").Value = 806.2
$ws.Range("AE5").Value = 734.2
$ws.Range("AG5").Value = 280.2
$ws.Range("AH5").Value = 665
$ws.Range("AI5").Value = 1160.8
$ws.Range("AJ5").Value = 5325
$ws.Range("AL5").Value = 958.6
$ws.Range("AM5").Value = 328
$ws.Range("AN5").Value = 84.09999999999999
$ws.Range("AO5").Value = 1

# Row 6
$ws.Range("C6").Value = 662425
$ws.Range("D6").Value = 1085.9
$ws.Range("E6").Value = 2935.8
$ws.Range("H6").Value = 333
$ws.Range("I6").Value = 1278
$ws.Range("J6").Value = 55081
$ws.Range("L6").Value = 1544.1
$ws.Range("M6").Value = 429
$ws.Range("N6").Value = 70.3
$ws.Range("O6").Value = 0.6
$ws.Range("P6").Value = 210798
$ws.Range("Q6").Value = 345.6
$ws.Range("R6").Value = 2192.3
$ws.Range("U6").Value = 50.5
$ws.Range("V6").Value = 221.2
$ws.Range("W6").Value = 48717
$ws.Range("Y6").Value = 483.5
$ws.Range("Z6").Value = 436
$ws.Range("AA6").Value = 71.5
$ws.Range("AC6").Value = 347639
$ws.Range("AD6").Value = 569.9
$ws.Range("AE6").Value = 742.4
$ws.Range("AG6").Value = 0.8
$ws.Range("AH6").Value = 330
$ws.Range("AI6").Value = 873
$ws.Range("AJ6").Value = 8295
$ws.Range("AL6").Value = 760.7
$ws.Range("AM6").Value = 457
$ws.Range("AN6").Value = 74.90000000000001
$ws.Range("AO6").Value = 0.4

# Row 7
$ws.Range("C7").Value = 419663
$ws.Range("D7").Value = 672.5
$ws.Range("E7").Value = 2419.3
$ws.Range("H7").Value = 65
$ws.Range("I7").Value = 620.8
$ws.Range("J7").Value = 50844
$ws.Range("L7").Value = 1128.1
$ws.Range("M7").Value = 372
$ws.Range("N7").Value = 59.6
$ws.Range("P7").Value = 103217
$ws.Range("Q7").Value = 165.4
$ws.Range("R7").Value = 624
$ws.Range("V7").Value = 86.5
$ws.Range("W7").Value = 10903
$ws.Range("Y7").Value = 276
$ws.Range("Z7").Value = 374
$ws.Range("AA7").Value = 59.9
$ws.Range("AC7").Value = 198030
$ws.Range("AD7").Value = 317.4
$ws.Range("AE7").Value = 540.6
$ws.Range("AH7").Value = 79
$ws.Range("AI7").Value = 423
$ws.Range("AJ7").Value = 4716
$ws.Range("AL7").Value = 521.1
$ws.Range("AM7").Value = 380
$ws.Range("AN7").Value = 60.9
$ws.Range("AO7").Value = -0.5

# Row 8
$ws.Range("C8").Value = 8970
$ws.Range("D8").Value = 560.6
$ws.Range("E8").Value = 1236.4
$ws.Range("I8").Value = 307.8
$ws.Range("J8").Value = 4576
$ws.Range("L8").Value = 1794
$ws.Range("P8").Value = 3339
$ws.Range("Q8").Value = 208.7
$ws.Range("R8").Value = 633.3
$ws.Range("V8").Value = 34.5
$ws.Range("W8").Value = 2528
$ws.Range("Y8").Value = 556.5
$ws.Range("Z8").Value = 6
$ws.Range("AA8").Value = 37.5
$ws.Range("AB8").Value = -0.8
$ws.Range("AC8").Value = 4401
$ws.Range("AD8").Value = 275.1
$ws.Range("AE8").Value = 419.4
$ws.Range("AH8").Value = 3
$ws.Range("AI8").Value = 374.8
$ws.Range("AJ8").Value = 1205
$ws.Range("AL8").Value = 550.1
$ws.Range("AM8").Value = 8
$ws.Range("AN8").Value = 50
$ws.Range("AO8").Value = -1.3

# Row 9
$ws.Range("AB9").Value = -2.3
$ws.Range("AC9").Value = 490
$ws.Range("AD9").Value = 490
$ws.Range("AF9").Value = 490
$ws.Range("AG9").Value = 490
$ws.Range("AH9").Value = 490
$ws.Range("AI9").Value = 490
$ws.Range("AJ9").Value = 490
$ws.Range("AL9").Value = 490

# Row 10
$ws.Range("C10").Value = 180386
$ws.Range("D10").Value = 1670.2
$ws.Range("E10").Value = 3362.9
$ws.Range("G10").Value = 16.2
$ws.Range("H10").Value = 589.5
$ws.Range("I10").Value = 1681.8
$ws.Range("J10").Value = 24791
$ws.Range("L10").Value = 2173.3
$ws.Range("P10").Value = 121802
$ws.Range("Q10").Value = 1127.8
$ws.Range("R10").Value = 3510.8
$ws.Range("T10").Value = 9.199999999999999
$ws.Range("U10").Value = 162
$ws.Range("V10").Value = 793.8
$ws.Range("W10").Value = 24554
$ws.Range("Y10").Value = 1433
$ws.Range("AB10").Value = 0.9
$ws.Range("AC10").Value = 72044
$ws.Range("AD10").Value = 667.1
$ws.Range("AE10").Value = 877.3
$ws.Range("AH10").Value = 399.5
$ws.Range("AI10").Value = 1024.8
$ws.Range("AJ10").Value = 6193
$ws.Range("AL10").Value = 923.6
$ws.Range("AM10").Value = 78
$ws.Range("AN10").Value = 72.2
$ws.Range("AO10").Value = 0.2

# Row 11
$ws.Range("C11").Value = 130730
$ws.Range("D11").Value = 514.7
$ws.Range("E11").Value = 924.7
$ws.Range("H11").Value = 35.5
$ws.Range("I11").Value = 699.2
$ws.Range("J11").Value = 5153
$ws.Range("L11").Value = 920.6
$ws.Range("P11").Value = 32742
$ws.Range("Q11").Value = 128.9
$ws.Range("R11").Value = 295.6
$ws.Range("U11").Value = 8
$ws.Range("V11").Value = 101
$ws.Range("W11").Value = 2073
$ws.Range("Y11").Value = 221.2
$ws.Range("Z11").Value = 148
$ws.Range("AA11").Value = 58.3
$ws.Range("AB11").Value = 0
$ws.Range("AC11").Value = 116988
$ws.Range("AD11").Value = 460.6
$ws.Range("AE11").Value = 614.4
$ws.Range("AH11").Value = 181
$ws.Range("AI11").Value = 725.5
$ws.Range("AJ11").Value = 3597
$ws.Range("AL11").Value = 722.1
$ws.Range("AM11").Value = 162
$ws.Range("AN11").Value = 63.8
$ws.Range("AO11").Value = -0.3

# Row 12
$ws.Range("C12").Value = 528852
$ws.Range("D12").Value = 440.3
$ws.Range("E12").Value = 1164
$ws.Range("I12").Value = 293
$ws.Range("J12").Value = 17617
$ws.Range("L12").Value = 942.7
$ws.Range("M12").Value = 561
$ws.Range("N12").Value = 46.7
$ws.Range("P12").Value = 195022
$ws.Range("Q12").Value = 162.4
$ws.Range("R12").Value = 761.1
$ws.Range("V12").Value = 48
$ws.Range("W12").Value = 13956
$ws.Range("Y12").Value = 338.6
$ws.Range("Z12").Value = 576
$ws.Range("AA12").Value = 48
$ws.Range("AB12").Value = -0.4
$ws.Range("AC12").Value = 340249
$ws.Range("AD12").Value = 283.3
$ws.Range("AE12").Value = 534.9
$ws.Range("AH12").Value = 15
$ws.Range("AI12").Value = 343
$ws.Range("AJ12").Value = 5226
$ws.Range("AL12").Value = 525.1
$ws.Range("AM12").Value = 648
$ws.Range("AN12").Value = 54
$ws.Range("AO12").Value = -1

# Row 13
$ws.Range("C13").Value = 134157
$ws.Range("D13").Value = 1081.9
$ws.Range("E13").Value = 4850
$ws.Range("H13").Value = 5
$ws.Range("I13").Value = 949.8
$ws.Range("J13").Value = 52263
$ws.Range("L13").Value = 2002.3
$ws.Range("M13").Value = 67
$ws.Range("N13").Value = 54
$ws.Range("P13").Value = 45586
$ws.Range("Q13").Value = 367.6
$ws.Range("R13").Value = 2247.9
$ws.Range("U13").Value = 3.5
$ws.Range("V13").Value = 117.5
$ws.Range("W13").Value = 24618
$ws.Range("Y13").Value = 701.3
$ws.Range("AC13").Value = 57064
$ws.Range("AD13").Value = 460.2
$ws.Range("AE13").Value = 697
$ws.Range("AH13").Value = 57.5
$ws.Range("AI13").Value = 758.2
$ws.Range("AJ13").Value = 3450
$ws.Range("AL13").Value = 792.6
$ws.Range("AM13").Value = 72
$ws.Range("AN13").Value = 58.1
$ws.Range("AO13").Value = -0.7

Write-Host "Updated 254 cells across 10 rows"
